$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update "Total Time" (column E) values for several tasks
$ws.Range("E2").Value = 2
$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 1

# These rows previously had no "Total Time" entry - add the value and match
# the centered-alignment formatting used by the other column-E cells.
$ws.Range("E7").Value = 1
$ws.Range("E7").HorizontalAlignment = $ws.Range("E6").HorizontalAlignment

$ws.Range("E8").Value = 1
$ws.Range("E8").HorizontalAlignment = $ws.Range("E6").HorizontalAlignment

$ws.Range("E9").Value = 1
$ws.Range("E9").HorizontalAlignment = $ws.Range("E6").HorizontalAlignment

$ws.Range("E10").Value = 2
$ws.Range("E10").HorizontalAlignment = $ws.Range("E6").HorizontalAlignment

$ws.Range("E11").Value = 1
$ws.Range("E11").HorizontalAlignment = $ws.Range("E6").HorizontalAlignment

$ws.Range("E12").Value = 2

# Update the active selection on the sheet
$ws.Range("K11").Select()
